$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the timestamp on the previously-last block of rows (450:463)
#    44232.34550189609 -> 44232.34550189815
for ($r = 450; $r -le 463; $r++) {
    $ws.Range("D$r").Value2 = 44232.34550189815
}

# 2) Append a brand-new refresh block: rows 464:477, same 14-row cycle of
#    Name/URL pairs used throughout the sheet, all stamped with the newest
#    check timestamp (44232.36656471538) and "Disponible" in column C.
$names = @("Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer", "Tomcat", "Shiny", "Github", "EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$stamp = 44232.36656471538
$startRow = 464

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("A$row").Value2 = $names[$i]

    $url = $urls[$i]
    $hashPos = $url.IndexOf("#")
    if ($hashPos -ge 0) {
        $address = $url.Substring(0, $hashPos)
        $subAddress = $url.Substring($hashPos + 1)
    } else {
        $address = $url
        $subAddress = ""
    }

    $ws.Range("B$row").Value2 = $url
    $target = $ws.Range("B$row")
    if ($subAddress -ne "") {
        $ws.Hyperlinks.Add($target, $address, $subAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($target, $address) | Out-Null
    }
    # Hyperlinks.Add re-applies font formatting on its own xf; reassert the
    # shared "Hyperlink" cell style afterwards so the cell keeps using the
    # workbook's existing Hyperlink style slot instead of the ad-hoc one.
    $ws.Range("B$row").Style = "Hyperlink"

    $ws.Range("C$row").Value2 = "Disponible"

    $ws.Range("D$row").Value2 = $stamp
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
